$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = ""

# Row 3
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""

# Row 4
$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""

# Row 5
$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("D5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""

# Row 6
$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("D6").Value = ""
$ws.Range("I6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("J6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("K6").Value = "7,97 TL - 15,96 TL - 199,41 TL"

# Row 7
$ws.Range("K7").Value = ""

# Row 8
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D8").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""

# Row 9
$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D9").Value = ""
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""

# Row 10
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("D10").Value = ""
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = ""

# Row 11
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("D11").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""

# Row 12
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("D12").Value = ""
$ws.Range("K12").Value = "WU: 12 USD–; Diğer: 404,16 TL–3.403,42 TL"

# Row 13
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = ""
$ws.Range("I13").Value = "Hesaba: Asgari 1 TL | Azami 7,97 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 69,62 TL"

# Row 14
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("J14").Value = "1.554,97 TL - 7.784 TL"
$ws.Range("K14").Value = "2.000 TL - 24.000 TL"
